$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.533.63'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '1.983.87'
$ws.Range('E3').Value = '  -3.60%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.07'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('E6').Value = '  -3.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.64'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -11.42%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.376'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.54'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -3.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +5.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.14'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +13.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.867'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -4.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.12'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -5.52%  '
$ws.Range('D16').Value = '2.273.92'
$ws.Range('E16').Value = '  -3.64%  '
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').Value = '1.987.27'
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').Value = '36.458.68'
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.17'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('D21').Value = '0.0₃0865'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.33'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.06'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.62'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  -3.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.25'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +5.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.22'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.88'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.126'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +9.31%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.92'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -6.93%  '
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.47'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -6.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.23'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +2.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.28'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -6.87%  '
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('E39').Value = '  -3.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.09'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0973'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -5.53%  '
$ws.Range('E43').Value = '  -3.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0214'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -2.71%  '
$ws.Range('E45').Value = '  -4.79%  '
$ws.Range('E46').Value = '  -4.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '92.89'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.64'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -3.87%  '
$ws.Range('D49').Value = '1.371.43'
$ws.Range('E49').Value = '  -3.39%  '
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.37'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -3.20%  '
